# daily auto push: 2026-02-15 22:42 UTC
# A new daily-ranking row (2026/02/16) is inserted right after the existing
# 2026/02/16 row (row 801), pushing rows 802:843 down to 803:844 and adding
# one brand-new trailing row at the end (old row 843's data moves to 844,
# new values are written into what is now row 843).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 802:843 down to 803:844 by inserting a new blank row at 802.
$ws.Rows.Item(802).Insert()

# Force column A to be treated as plain text so the "yyyy/mm/dd"-looking
# string isn't auto-converted into a date serial number/date format (the
# rest of the sheet stores these as plain text too).
$ws.Cells.Item(802, 1).NumberFormat = "@"
$ws.Cells.Item(802, 1).Value = "2026/02/16"
# Drop back to the default "Normal" style so no stray per-cell formatting
# is left behind, matching the rest of the data rows (which carry no
# explicit style).
$ws.Cells.Item(802, 1).Style = "Normal"

$ws.Cells.Item(802, 2).Value = "月"
$ws.Cells.Item(802, 3).Value = 3
$ws.Cells.Item(802, 4).Value = 36
